# Updating for autumn 2024
# -------------------------------------------------------------------------
# This workbook tracked a single year's schedule (Schedule_date /
# module_due_dates). For Autumn 2024 we keep last year's sheets around
# (renamed *_2023) and turn the original two sheets into the live 2024
# versions with refreshed dates.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------
# 1) Duplicate the two date-driven sheets BEFORE editing them, so the
#    copies keep the original (2023) data, then rename everything.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Schedule_date").Copy($missing, $wb.Worksheets.Item("Schedule_date"))
$wb.Worksheets.Item("module_due_dates").Copy($missing, $wb.Worksheets.Item("module_due_dates"))

$wb.Worksheets.Item("Schedule_date").Name = "Schedule_date_2024"
$wb.Worksheets.Item("module_due_dates").Name = "module_due_dates_2024"
$wb.Worksheets.Item("Schedule_date (2)").Name = "Schedule_date_2023"
$wb.Worksheets.Item("module_due_dates (2)").Name = "module_due_dates_2023"

# Put module_due_dates_2024 right after Schedule_date_2024 so the tab
# order reads: Schedule, Schedule_date_2024, module_due_dates_2024,
# Schedule_date_2023, module_due_dates_2023
$wb.Worksheets.Item("module_due_dates_2024").Move($missing, $wb.Worksheets.Item("Schedule_date_2024"))

# ---------------------------------------------------------------------
# 2) Schedule_date_2024 - shift every class date forward to Autumn 2024
#    (the weekly Tuesday cadence just moves a year later), and fix up
#    the Thanksgiving break which now lands on row 16 instead of row 15.
# ---------------------------------------------------------------------
$sched2024 = $wb.Worksheets.Item("Schedule_date_2024")

$dates2024 = @(45524, 45531, 45538, 45545, 45552, 45559, 45566, 45573, 45580, 45587, 45594, 45601, 45608, 45615, 45622, 45629)
for ($i = 0; $i -lt $dates2024.Length; $i++) {
    $row = $i + 2
    $sched2024.Range("A$row").Value = $dates2024[$i]
}

# Row 15 becomes the regular "4: Putting it together" / extension-packages week
$sched2024.Range("B15").Value = "4: Putting it together"
$sched2024.Range("C15").Value = "ggplot extension packages and complexheatmap"
# Row 16 becomes the Thanksgiving break week
$sched2024.Range("B16").Value = "No class, Thanksgiving"
$sched2024.Range("C16").Value = "Relaxing and eating"

$sched2024.Range("A1:C17").Select() | Out-Null
$sched2024.Range("A12").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) module_due_dates_2024 - replace the due-date table with the 2024
#    version: a new "Reflections" row plus real 2024 due dates (written
#    as text, e.g. "Monday, August 26, 2024").
# ---------------------------------------------------------------------
$due2024 = $wb.Worksheets.Item("module_due_dates_2024")

$due2024.Cells.Clear()

$due2024.Range("A1").Value = "Assignment"
$due2024.Range("B1").Value = "Due Date"

$due2024.Range("A2").Value = "Reflections"
$due2024.Range("B2").Value = "1 week after each class"

$due2024.Range("A3").Value = "Module 1: Good and bad visualizations"
$due2024.Range("B3").Value = "Monday, August 26, 2024"

$due2024.Range("A4").Value = "Module 2: Coding Fundamentals"
$due2024.Range("B4").Value = "Tuesday, October 1, 2024"

$due2024.Range("A5").Value = "Module 3: Data Exploration"
$due2024.Range("B5").Value = "Tuesday, October 29, 2024"

$due2024.Range("A6").Value = "Module 4: Putting it together"
$due2024.Range("B6").Value = "Tuesday, December 3, 2024"

$due2024.Range("A7").Value = "Capstone plan"
$due2024.Range("B7").Value = "Tuesday, November 5, 2024"

$due2024.Range("A8").Value = "Capstone"
$due2024.Range("B8").Value = "Friday, December 6, 2024"

# (the host engine snaps ColumnWidth to whole pixels, so these inputs are
# pre-compensated to land as close as possible to the real bestFit widths)
$due2024.Columns.Item(1).ColumnWidth = 25.5
$due2024.Columns.Item(2).ColumnWidth = 23.3
$due2024.Columns.Item(3).ColumnWidth = 21.3

# ---------------------------------------------------------------------
# 4) Leave the *_2023 archive sheets' selections roughly where the
#    original sheets had them.
# ---------------------------------------------------------------------
$sched2023 = $wb.Worksheets.Item("Schedule_date_2023")
$sched2023.Activate()
$sched2023.Range("E32").Select() | Out-Null

$due2023 = $wb.Worksheets.Item("module_due_dates_2023")
$due2023.Activate()
$due2023.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) Finally, land on module_due_dates_2024 - the sheet this edit is
#    really about - matching the workbook's saved active tab.
# ---------------------------------------------------------------------
$due2024.Activate()
$due2024.Range("B5").Select() | Out-Null

